# "preping for resubmission to nat micro"
# The active sheet (ForAnalysis) gets its view scrolled/re-selected and
# column A widened so the strain names are fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ForAnalysis")
$ws.Activate()

# Widen column A (stored OOXML width 30.5 characters).
$ws.Columns.Item(1).ColumnWidth = 29.7

# Scroll the view down one row and move the selection to B26, matching
# where the user left off reviewing the data before resubmission.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("B26").Select()
